# Update row-2 metric values on each model sheet so that all sheets share
# the same train_test_split (note column B, the "Arreglo aleatorio óptimo",
# now equals 16 everywhere).

$wb = $excel.ActiveWorkbook

$updates = @{
    "arbolu"   = @{ "A2"=14;  "B2"=16; "C2"=0.4705882352941176;  "D2"=0;                    "E2"=0.4705882352941176;  "F2"=0.4152661064425771;  "G2"=0.4705882352941176; "H2"=0.4410427807486632; "I2"=0.6666666666666666 }
    "bosqueu"  = @{ "A2"=106; "B2"=16; "C2"=0.6470588235294118;  "D2"=0;                    "E2"=0.6764705882352942;  "F2"=0.6176470588235294;  "G2"=0.6764705882352942; "H2"=0.6211636828644501; "I2"=0.8300173010380623 }
    "knnu"     = @{ "A2"=23;  "B2"=16; "C2"=0.5882352941176471;  "D2"=0;                    "E2"=0.5882352941176471;  "F2"=0.6158645276292335;  "G2"=0.5882352941176471; "H2"=0.461846976552859;  "I2"=0.7869088811995386 }

    "arbolts"  = @{ "A2"=3;   "B2"=16; "C2"=0.6013077211281678;  "D2"=0.5281202066941012;  "E2"=0.7267187947852327;  "F2"=-0.3424411177081306 }
    "bosquets" = @{ "A2"=117; "B2"=16; "C2"=0.5506265049906924;  "D2"=0.4980777357456132;  "E2"=0.7057462261646272;  "F2"=-0.2660754574519941 }
    "knnts"    = @{ "A2"=5;   "B2"=16; "C2"=0.4954067924349528;  "D2"=0.5209749223687926;  "E2"=0.7217859255823659;  "F2"=-0.3242783521967445 }

    "arboltd"  = @{ "A2"=16;  "B2"=16; "C2"=0.7208162461912868;  "D2"=0.832133534832235;   "E2"=0.9122135357646448;  "F2"=-0.05700681345855352 }
    "bosquetd" = @{ "A2"=113; "B2"=16; "C2"=0.5810498849275845;  "D2"=0.610917145371212;   "E2"=0.7816118892207385;  "F2"=0.2239904317205436 }
    "knntd"    = @{ "A2"=3;   "B2"=16; "C2"=0.6918548791568156;  "D2"=0.8446574183655695;  "E2"=0.9190524568084074;  "F2"=-0.07291511383532523 }

    "arbolcc"  = @{ "A2"=5;   "B2"=16; "C2"=0.6942176303905359;  "D2"=0.8853505299486987;  "E2"=0.9409306722329221;  "F2"=-0.1108380796448769 }
    "bosquecc" = @{ "A2"=110; "B2"=16; "C2"=0.4532064780878013;  "D2"=0.440065233326777;   "E2"=0.6633741277188741;  "F2"=0.310624504167525 }
    "knncc"    = @{ "A2"=5;   "B2"=16; "C2"=0.5118031337185958;  "D2"=0.5518818102584471;  "E2"=0.7428874815599245;  "F2"=0.3075608930465602 }

    "arbolpp"  = @{ "A2"=22;  "B2"=16; "C2"=0.4826264106791714;  "D2"=1.037497801426076;   "E2"=1.018576360135103;   "F2"=-0.6805716499436505 }
    "bosquepp" = @{ "A2"=118; "B2"=16; "C2"=0.477622994526035;   "D2"=0.3783366230554467;  "E2"=0.6150907437569246;  "F2"=0.3871584094265622 }
    "knnpp"    = @{ "A2"=9;   "B2"=16; "C2"=0.6247532853511617;  "D2"=0.5509014783289405;  "E2"=0.742227376434567;   "F2"=-0.02460949500583753 }
}

foreach ($sheetName in $updates.Keys) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $cellValues = $updates[$sheetName]
    foreach ($addr in $cellValues.Keys) {
        $sheet.Range($addr).Value = $cellValues[$addr]
    }
}
